$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab
$ws.Name = "SA"

# 2) Tiny floating point precision corrections on existing rows
$ws.Range("C13").Value = 0.995526852596242
$ws.Range("G13").Value = 0.995526852596242
$ws.Range("M13").Value = 0.9954987230375346
$ws.Range("P13").Value = 0.9920590896461766

$ws.Range("C15").Value = 0.9304766936225563
$ws.Range("E15").Value = 0.9555256776006366
$ws.Range("G15").Value = 0.9304766936225563
$ws.Range("O15").Value = 0.9988553507919351

# 3) Append new row 16 with data (row 14 / HexGrid-60degTilt5degRes series)
# Copy formatting from row 15's A cell (bold border style) before writing the
# new value, since direct .Style assignment does not propagate reliably.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.9783622755493311
$ws.Range("D16").Value = 0.7742431438492958
$ws.Range("E16").Value = 0.8977956654838386
$ws.Range("F16").Value = 0.8651770259480333
$ws.Range("G16").Value = 0.9783622755493311
$ws.Range("H16").Value = 0.7742431438492958
$ws.Range("I16").Value = 1.17407315656008
$ws.Range("J16").Value = 0.8497568785102666
$ws.Range("K16").Value = 1.078894377174901
$ws.Range("L16").Value = 0.8964656254791129
$ws.Range("M16").Value = 0.9783622755493311
$ws.Range("N16").Value = 0.8360194046665672
$ws.Range("O16").Value = 0.8788945277076247
$ws.Range("P16").Value = 0.9393460185693574
